$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Finalize existing row 54 (IN PROGRESS -> DONE) ---
$ws.Range("H54").Value = "DONE"
$ws.Range("I54").Value = 42866.613298611112
$ws.Range("J54").Value = "0.04987508 USDT (0.15%) "
$ws.Range("K54").Value = "     ~11%"
$ws.Range("K54").Font.Color = 5287936
$ws.Range("L54").Value = " 2 day"

# --- Append new row 55 (new Buy transaction, IN PROGRESS) ---
$ws.Range("A55").Value = 42867.279143518521
$ws.Range("B55").Value = "            Buy"
$chars = $ws.Range("B55").Characters(13, 3)
$chars.Font.Color = 5287936
$ws.Range("C55").Value = "        XRP"
$ws.Range("D55").Value = "        0.192
"
$ws.Range("E55").Value = "         0.175USDT"
$ws.Range("F55").Value = "         189 XRP"
$ws.Range("G55").Value = " XRP/USDT0000006"
$ws.Range("H55").Value = "IN PROGRESS"
$ws.Range("K55").Value = "     "

# --- Scroll position update ---
$ws.Application.ActiveWindow.ScrollRow = 43
